# Fill in the remaining survey responses (column B) for rows 208-301 on the
# "Treinamento" sheet, then leave the view scrolled/selected at the bottom,
# mirroring the state Excel saves after the user finishes entering data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

$values = @(0,1,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,1,1,1,1,0,0,1,1,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,1,0,0,1,0,0,0,0,0,0,0,1,0,0,0,0)

$startRow = 208
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Activate the sheet and set the view/selection the way it ended up after
# the last edit (scrolled to show the final rows, cursor on the first empty
# cell below the filled data).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 293
$ws.Range("B302").Select()
